{"js": "// 1. Bump the GroupDocs.Assembly evaluation-watermark version: 25.6 -> 25.12\nconst results = context.document.body.search(\"GroupDocs.Assembly 25.6.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nresults.items.forEach((r) => r.insertText(\"GroupDocs.Assembly 25.12.\", Word.InsertLocation.replace));\nawait context.sync();\n\n// 2. Add the (previously missing) built-in \"Hyperlink\" character style, based on\n//    \"Default Paragraph Font\", blue + single underline - matches what Word itself\n//    mints the first time a hyperlink is inserted into the document.\ncontext.document.addStyle(\"Hyperlink\", Word.StyleType.character);\nawait context.sync();\n\nconst style = context.document.getStyles().getByNameOrNullObject(\"Hyperlink\");\nstyle.baseStyle = \"DefaultParagraphFont\";\nstyle.font.color = \"#0000FF\";\nstyle.font.underline = Word.UnderlineType.single;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Bump the GroupDocs.Assembly evaluation-watermark version: 25.6 -> 25.12\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"GroupDocs.Assembly 25.6.\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"GroupDocs.Assembly 25.12.\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 2. Add the (previously missing) built-in \"Hyperlink\" character style, based on\n#    \"Default Paragraph Font\", blue + single underline - matches what Word itself\n#    mints the first time a hyperlink is inserted into the document.\n$style = $d.Styles.Add(\"Hyperlink\", 2)\n$style.BaseStyle = \"DefaultParagraphFont\"\n$style.Font.Color = 16711680\n$style.Font.Underline = 1\n"}
